# Generate Report for Handoff
#
# This mirrors the localization-status report moving from "In Translation"
# to "Ready for handoff": the status text and the two "last generated"
# timestamps are refreshed, and the now-wider status column is resized to
# fit the new label on the Overview/zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------
# Overview!E2 (zh-cn status) and Overview!F2 (de-de status)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# zh-cn!C2 / de-de!C2 "Status" column
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Latest HO Xliff Generate Date / Latest Handoff Datetime ----------
# Overview!G2 and de-de!H2 shared the same timestamp string.
$wsOverview.Range("G2").Value = "2016-08-21 05:04:09"
$wsDeDe.Range("H2").Value = "2016-08-21 05:04:09"
# zh-cn!H2 had its own timestamp string.
$wsZhCn.Range("H2").Value = "2016-08-21 05:04:03"

# --- Widen the Status / status columns to fit "Ready for handoff" -----
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # column F
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # column C
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # column C
